$d = $word.ActiveDocument

# --- 1) "Expected mark" cell for Pedro Martinho (student 78141, table row 2)
#        gets a new run with text "                 7" (leading spaces + 7) ---
$markCellMartinho = $d.Tables(1).Rows(2).Cells(4)
$paraMartinho = $markCellMartinho.Range.Paragraphs(1)
$rngMartinho = $paraMartinho.Range
$rngMartinho.InsertAfter("                 7")
$rngMartinho.Font.Size = 8
$rngMartinho.Font.SizeBi = 8

# --- 2) "Expected mark" cell for Rui Figueiredo (student 78247, table row 3)
#        "6.5" becomes "6." + "8" (two separate runs) i.e. value 6.8 ---
$markCellFigueiredo = $d.Tables(1).Rows(3).Cells(4)
$paraFigueiredo = $markCellFigueiredo.Range.Paragraphs(1)
$fullFigueiredo = $paraFigueiredo.Range
$lastDigit = $d.Range($fullFigueiredo.Start + 2, $fullFigueiredo.Start + 3)
$lastDigit.Text = "8"
# force the replaced character into its own run (same resulting formatting)
$lastDigit2 = $d.Range($fullFigueiredo.Start + 2, $fullFigueiredo.Start + 3)
$lastDigit2.Bold = 1
$lastDigit2.Bold = 0

# --- 3) move the "_GoBack" bookmark from the paragraph right after the table
#        into the first (empty) paragraph of the "Expected mark" cell for
#        Alexandre Candeias (student 78599, table row 4) ---
$markCellCandeias = $d.Tables(1).Rows(4).Cells(4)
$paraCandeias = $markCellCandeias.Range.Paragraphs(1)
$d.Bookmarks.Add("_GoBack", $paraCandeias.Range)
